$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("On-Site")

# Fill in the new "Umbilical Cable" sub-system row (row 5), matching the
# existing data rows above it.
$ws.Range("A5").Value = "Umbilical Cable"
$ws.Range("B5").Value = 30000
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 6

# Give the new row the same cell style used by the other sub-system rows.
$ws.Range("A2:A5").Style = "Normal 2"

# Leave the selection on the newly added cell, as in the source workbook.
$ws.Range("A5").Select()
